# Fixed Bento 80 Test scripts
# Append an "order By ... ASC LIMIT 100" clause to the end of each tab's
# Neo4j query (column B) on the "startup" sheet, for CasesTab (row 2),
# SamplesTab (row 3) and FilesTab (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- FilesTab (row 4) ---
$b4 = $ws.Range("B4").Value2
$ws.Range("B4").Value = $b4 + "`norder By f.file_name ASC LIMIT 100"

# --- SamplesTab (row 3) ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- CasesTab (row 2) ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value = $b2 + "`norder By ss.study_subject_id ASC LIMIT 100"

# Restore the selection/view state referenced by the edit (active cell C3,
# no pinned top-left cell).
$ws.Range("C3").Select()
